$d = $word.ActiveDocument

# Minimal single-part WordprocessingML package wrapper used by every
# InsertXML call below (Range.InsertXML needs a full package, not a
# bare <w:document> fragment).
$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Find-ParagraphByPrefix($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ($prefix + "*")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) The paragraph "IdInventario quedo como pk ..." currently holds its
#    text split across two runs ("... SerieKey en " / "Inventario")
#    with the _GoBack bookmark sandwiched in between. Collapse it to a
#    single merged run and drop the bookmark here (it gets re-created
#    further down, after the new content we are about to insert).
# ---------------------------------------------------------------------
$target = Find-ParagraphByPrefix "IdInventario quedo como pk de Inventario"
if ($null -eq $target) {
    throw "Could not find the 'IdInventario quedo como pk de Inventario' paragraph"
}

$mergedRunXml = '<w:body><w:p><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>IdInventario quedo como pk de Inventario y hay un índice unique entre IdBienEspecif y SerieKey en Inventario</w:t></w:r></w:p></w:body>'
$range = $d.Range($target.Range.Start, $target.Range.End - 1)
$range.InsertXML($pkgHeader + $mergedRunXml + $pkgFooter)

# Re-resolve the paragraph: its content changed so re-fetch a fresh
# reference before using it as an insertion anchor.
$target = Find-ParagraphByPrefix "IdInventario quedo como pk de Inventario"

# ---------------------------------------------------------------------
# 2) Insert the new block of annotation paragraphs right after it.
# ---------------------------------------------------------------------
$sz = '<w:sz w:val="16"/><w:szCs w:val="16"/>'
$szHl = $sz + '<w:highlight w:val="yellow"/>'
$numPr = '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>'

$newParas = ''

# "Elimine columna IdBien ..." - numbered, not highlighted
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/>' + $numPr + '<w:rPr>' + $sz + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $sz + '</w:rPr><w:t>Elimine columna IdBien en AsigDetalle, porque la pk de inventario ya no contiene a IdBien</w:t></w:r></w:p>'

# "Quitar el campo Observacion ..." - numbered, highlighted
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/>' + $numPr + '<w:rPr>' + $szHl + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $szHl + '</w:rPr><w:t>Quitar el campo Observacion de AsigDetalle (poner solo uno en asignación quizas)</w:t></w:r></w:p>'

# "//FALTA QUE Ponga en estado ..." - numbered, highlighted
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/>' + $numPr + '<w:rPr>' + $szHl + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $szHl + '</w:rPr><w:t>//FALTA QUE Ponga en estado &quot;Compradas&quot; las cosas de un SolicDetalle</w:t></w:r></w:p>'

# "//Para hacerlo se puede consultar ..." - not numbered, highlighted
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr>' + $szHl + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $szHl + '</w:rPr><w:t>//Para hacerlo se puede consultar de nuevo la cantidad comprada y comparar con la cantidad solicitada</w:t></w:r></w:p>'

# "//Todo dentro del negocio o dal" + ".. El código ..." - not numbered, highlighted, two runs
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr>' + $szHl + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $szHl + '</w:rPr><w:t>//Todo dentro del negocio o dal</w:t></w:r>' +
    '<w:r><w:rPr>' + $szHl + '</w:rPr><w:t>.. El código donde empzaria esto esa en frmregistrar línea 66</w:t></w:r></w:p>'

# "Al registrar bienes, ..." - numbered, highlighted, has the (now moved) lastRenderedPageBreak
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/>' + $numPr + '<w:rPr>' + $szHl + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $szHl + '</w:rPr><w:lastRenderedPageBreak/><w:t>Al registrar bienes, solo me tiene que permitir registrar los de un mismo SolicDetalle (sino el proveedor no concuerda con todos los datalles)</w:t></w:r></w:p>'

# "Seguir en DALAsignacion ..." - numbered, highlighted; _GoBack bookmark now lives here, after the run
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/>' + $numPr + '<w:rPr>' + $szHl + '</w:rPr></w:pPr>' +
    '<w:r><w:rPr>' + $szHl + '</w:rPr><w:t>Seguir en DALAsignacion línea 58, lo de poner como finalizado un solicDetalle</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# trailing empty paragraph - not numbered, highlighted, no runs
$newParas += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr>' + $szHl + '</w:rPr></w:pPr></w:p>'

$insertXml = $pkgHeader + '<w:body>' + $newParas + '</w:body>' + $pkgFooter

$insertionPoint = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$insertionPoint.InsertXML($insertXml)

# ---------------------------------------------------------------------
# 3) The page-break hint that used to sit in front of the "When
#    clicking " run (further down the doc) now belongs on the new "Al
#    registrar bienes" paragraph instead, so drop it from here. Rebuild
#    that paragraph's run content byte-for-byte (including its
#    non-breaking spaces) minus the <w:lastRenderedPageBreak/> marker.
# ---------------------------------------------------------------------
$wc = Find-ParagraphByPrefix "When clicking"
if ($null -ne $wc) {
    $wcRuns = '<w:r w:rsidRPr="00143796"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="242729"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:lang w:val="en-US" w:eastAsia="es-AR"/></w:rPr><w:t>When clicking </w:t></w:r><w:r w:rsidRPr="00143796"><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="242729"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:shd w:val="clear" w:color="auto" w:fill="EFF0F1"/><w:lang w:val="en-US" w:eastAsia="es-AR"/></w:rPr><w:t>Button1</w:t></w:r><w:r w:rsidRPr="00143796"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="242729"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:lang w:val="en-US" w:eastAsia="es-AR"/></w:rPr><w:t> on the UserControl, i''ll fire </w:t></w:r><w:r w:rsidRPr="00143796"><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="242729"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:shd w:val="clear" w:color="auto" w:fill="EFF0F1"/><w:lang w:val="en-US" w:eastAsia="es-AR"/></w:rPr><w:t>Button1_Click</w:t></w:r><w:r w:rsidRPr="00143796"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="242729"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:lang w:val="en-US" w:eastAsia="es-AR"/></w:rPr><w:t> which triggers </w:t></w:r><w:r w:rsidRPr="00143796"><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="242729"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:shd w:val="clear" w:color="auto" w:fill="EFF0F1"/><w:lang w:val="en-US" w:eastAsia="es-AR"/></w:rPr><w:t>UserControl_ButtonClick</w:t></w:r><w:r w:rsidRPr="00143796"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="242729"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:lang w:val="en-US" w:eastAsia="es-AR"/></w:rPr><w:t> on the form:</w:t></w:r>'
    $wcXml = $pkgHeader + '<w:body><w:p>' + $wcRuns + '</w:p></w:body>' + $pkgFooter
    $wcRange = $d.Range($wc.Range.Start, $wc.Range.End - 1)
    $wcRange.InsertXML($wcXml)
}

Write-Output "ok"
